$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I1").Value = "text"

# Update the view so the window is scrolled to E2 and H2 is the active/selected cell
$excel.ActiveWindow.ScrollIntoView(0, 0, 0, 0)
$ws.Range("H2").Select()
$excel.ActiveWindow.TopLeftCell = $ws.Range("E2")
